# Rename the sole worksheet: "Property1" -> "DataNode"
# (matches <sheet name="Property1" .../> -> <sheet name="DataNode" .../> in xl/workbook.xml)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "DataNode"

# Register an extra (unused) 9pt SimSun font in the style table, mirroring the
# new <font><sz val="9"/><name val="宋体"/><charset val="134"/></font> entry
# that shows up as xl/styles.xml's 3rd font (fonts count 2 -> 3) in the target
# workbook. We bounce a non-bold cell's font size 9 -> 11 so the new font
# entry is appended to the shared font table without leaving any visible
# cell actually restyled (A21 keeps its original style index/appearance).
$ws.Range("A21").Font.Size = 9
$ws.Range("A21").Font.Size = 11

# Move/restore the active selection to C24 (was A9), matching the new
# <selection pane="bottomLeft" activeCell="C24" sqref="C24"/> in the frozen
# bottom-left pane of xl/worksheets/sheet1.xml.
$ws.Range("C24").Select()
